$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.857.32'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.499.66'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.37'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.11'
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.40'
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.347'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').Value = '2.941.64'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').Value = '58.793.63'
$ws.Range('E14').Value = '  +0.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.73'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000138'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = '2.502.03'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.03'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.47'
$ws.Range('E20').Value = '  -0.59%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.93'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('E23').Value = '  +1.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.419'
$ws.Range('E24').Value = '  +0.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('E27').Value = '  -0.58%  '
$ws.Range('D28').Value = '0.0₃0760'
$ws.Range('E28').Value = '  -1.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.82'
$ws.Range('E29').Value = '  +1.81%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.45'
$ws.Range('E30').Value = '  -3.29%  '
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('E32').Value = '  +0.08%  '
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.31'
$ws.Range('E34').Value = '  -0.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.36'
$ws.Range('E35').Value = '  -1.00%  '
$ws.Range('E36').Value = '  -1.02%  '
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '281.66'
$ws.Range('E41').Value = '  +0.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '10.92'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.96'
$ws.Range('E44').Value = '  -5.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '128.50'
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('E48').Value = '  -1.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.20'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '1.748.40'
$ws.Range('E51').Value = '  -0.51%  '
